$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert two blank rows above the current header row (row 1) ---
# Before: row1=header, rows2-5=data(Cell1..Cell4, "implementation" values)
# After two inserts: row1=blank, row2=blank, row3=header, rows4-7=data
$ws.Rows("1:1").Insert()
$ws.Rows("1:1").Insert()

# --- Duplicate header + original ("Implementation") data block further down ---
# Source A3:F7 (header + 4 data rows) -> destination starting at A11 (rows 11-15)
$ws.Range("A3:F7").Copy($ws.Range("A11"))

# --- Overwrite rows 4-7 with the new ("Paper") measurement values ---
$ws.Range("B4").Value = 2.506
$ws.Range("C4").Value = 0.88
$ws.Range("D4").Value = 1.101
$ws.Range("E4").Value = 95
$ws.Range("F4").Value = 0.975

$ws.Range("B5").Value = 1.951
$ws.Range("C5").Value = 0.935
$ws.Range("D5").Value = 1.044
$ws.Range("E5").Value = 75
$ws.Range("F5").Value = 0.672

$ws.Range("B6").Value = 1.655
$ws.Range("C6").Value = 0.786
$ws.Range("D6").Value = 0.908
$ws.Range("E6").Value = 76.316
$ws.Range("F6").Value = 0.739

$ws.Range("B7").Value = 2.809
$ws.Range("C7").Value = 1.462
$ws.Range("D7").Value = 1.62
$ws.Range("E7").Value = 46.875
$ws.Range("F7").Value = 0.836

# --- Add section title rows: "Implementation" (row 10) and "Paper" (row 2) ---
# NOTE: "Implementation" is entered before "Paper" so the shared-string table
# ends up ordered the same way the authored workbook has it (Implementation=10, Paper=11).
$ws.Range("A10").Value = "Implementation"
$ws.Range("A10:F10").Merge()
$ws.Range("A10:F10").Font.Bold = $true
$ws.Range("A10:F10").Font.Family = 3
$ws.Range("A10:F10").HorizontalAlignment = -4108
$ws.Range("A10:F10").VerticalAlignment = -4108
$ws.Range("A10:F10").Borders.Item(9).LineStyle = 1

$ws.Range("A2").Value = "Paper"
$ws.Range("A2:F2").Merge()
$ws.Range("A2:F2").Font.Bold = $true
$ws.Range("A2:F2").Font.Family = 3
$ws.Range("A2:F2").HorizontalAlignment = -4108
$ws.Range("A2:F2").VerticalAlignment = -4108
$ws.Range("A2:F2").Borders.Item(9).LineStyle = 1

# --- Selection cursor, matching the authored workbook state ---
$ws.Range("I9").Select()
